# Commit: "Add files via upload"
# Consolidate the COP 2025 / COP 2027 / COP 2028 campaign rows into "COP 2026",
# extend the long Details note in row 5, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Campaign group) consolidation -> "COP 2026"
$ws.Range("A7").Value  = "COP 2026"   # was "COP 2027"
$ws.Range("A9").Value  = "COP 2026"   # was "COP 2025"
$ws.Range("A10").Value = "COP 2026"   # was "COP 2025"
$ws.Range("A11").Value = "COP 2026"   # was "COP 2025"
$ws.Range("A12").Value = "COP 2026"   # was "COP 2028"

# Extend the Details text on row 5 (Digital Display Phase 1)
$current = $ws.Range("E5").Value()
$ws.Range("E5").Value = $current + ". Normally this wouldn't be this long"

# Update the selected range shown when the file was last saved
$ws.Range("A3:A12").Select()
